# daily auto push: 2026-01-06 18:47 UTC
# A new reading for 2026/01/06 is inserted in its correct chronological
# position (right after the existing 2026/01/06 rows, before the
# 2026/12/29 block), pushing every following row down by one, and a new
# trailing reading is appended for the last existing date (2027/01/05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 579; this shifts rows 579..620 down to 580..621
# and automatically extends the used range / dimension to D621.
$ws.Rows.Item(579).EntireRow.Insert()

# Fill in the newly inserted row with the new data point.
# The leading apostrophe forces the date-like string to be stored as
# plain text (matching the rest of column A) instead of being
# auto-converted into a date serial number by Excel's input parsing.
$ws.Range("A579").Value = "'2026/01/06"
$ws.Range("B579").Value = "火"
$ws.Range("C579").Value = 22
$ws.Range("D579").Value = 201
